$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "out_vars" (sheet1): add June 8 (serial 43990) row of national
# summary numbers; resize column A; move selection.
# ---------------------------------------------------------------------
$wsOut = $wb.Worksheets.Item("out_vars")
$wsOut.Activate()

# Clone the formatting of the prior data row down into the new row so the
# new cells pick up matching (wrap-text) styling, then overwrite values.
$wsOut.Range("A8:J8").Copy()
$wsOut.Range("A9:J9").PasteSpecial(-4122)

$wsOut.Range("A9").Value = 43990
$wsOut.Range("B9").Value = 120102
$wsOut.Range("C9").Value = 177875
$wsOut.Range("D9").Value = 46398
$wsOut.Range("E9").Value = 14053
$wsOut.Range("F9").Value = 33.315015570098751
$wsOut.Range("G9").Value = 40012
$wsOut.Range("H9").Value = 3733
$wsOut.Range("I9").Value = 3891
$wsOut.Range("J9").Value = 344375

# Column A narrows from "bestFit" 23.33 down to 10.5 characters.
$wsOut.Columns("A").ColumnWidth = 9.6

$wsOut.Range("E23").Select()

# ---------------------------------------------------------------------
# Sheet "dates_dx" (sheet2): row 9 (already present, blank) gets values.
# ---------------------------------------------------------------------
$wsDx = $wb.Worksheets.Item("dates_dx")
$wsDx.Activate()

$wsDx.Range("A9").Value = 43990
$wsDx.Range("B9").Value = 0
$wsDx.Range("C9").Value = 1
$wsDx.Range("D9").Value = 1
$wsDx.Range("E9").Value = 1
$wsDx.Range("F9").Value = 0
$wsDx.Range("G9").Value = 0
$wsDx.Range("H9").Value = 0
$wsDx.Range("I9").Value = 4

$wsDx.Range("G11").Select()

# ---------------------------------------------------------------------
# Sheet "dates_sx" (sheet3): new row 9 of data + a fresh, still-empty
# row 10 underneath (date-formatted like column A above it).
# ---------------------------------------------------------------------
$wsSx = $wb.Worksheets.Item("dates_sx")
$wsSx.Activate()

$wsSx.Range("A2").Copy()
$wsSx.Range("A9:A10").PasteSpecial(-4122)

$wsSx.Range("A9").Value = 43990
$wsSx.Range("B9").Value = 0
$wsSx.Range("C9").Value = 1
$wsSx.Range("D9").Value = 0
$wsSx.Range("E9").Value = 1
$wsSx.Range("F9").Value = 1
$wsSx.Range("G9").Value = 1
$wsSx.Range("H9").Value = 0
$wsSx.Range("I9").Value = 1
$wsSx.Range("J9").Value = 1
$wsSx.Range("K9").Value = 0
$wsSx.Range("L9").Value = 0

$wsSx.Range("D19").Select()

# ---------------------------------------------------------------------
# Sheet "dates_deaths" (sheet4): new row 9 of data.
# ---------------------------------------------------------------------
$wsDeaths = $wb.Worksheets.Item("dates_deaths")
$wsDeaths.Activate()

$wsDeaths.Range("A2").Copy()
$wsDeaths.Range("A9").PasteSpecial(-4122)

$wsDeaths.Range("A9").Value = 43990
$wsDeaths.Range("B9").Value = 1
$wsDeaths.Range("C9").Value = 0
$wsDeaths.Range("D9").Value = 2
$wsDeaths.Range("E9").Value = 1
$wsDeaths.Range("F9").Value = 1
$wsDeaths.Range("G9").Value = 2
$wsDeaths.Range("H9").Value = 2

$wsDeaths.Range("I9").Select()

# ---------------------------------------------------------------------
# Sheet "control_obs" (sheet5): new "2020-06-08" column (I) filled in
# for every existing metric row, plus the grand-total formula extended.
# ---------------------------------------------------------------------
$wsCtrl = $wb.Worksheets.Item("control_obs")
$wsCtrl.Activate()

$wsCtrl.Range("I1").Value = 43990

$wsCtrl.Range("I2").Value = 3273
$wsCtrl.Range("I3").Value = 3093
$wsCtrl.Range("I4").Value = 3093
$wsCtrl.Range("I5").Value = 3093
$wsCtrl.Range("I6").Value = 3093
$wsCtrl.Range("I7").Value = 2363
$wsCtrl.Range("I8").Value = 4910

$wsCtrl.Range("I10").Value = 148
$wsCtrl.Range("I11").Value = 148
$wsCtrl.Range("I12").Value = 148
$wsCtrl.Range("I13").Value = 148
$wsCtrl.Range("I14").Value = 148
$wsCtrl.Range("I15").Value = 125
$wsCtrl.Range("I16").Value = 160

$wsCtrl.Range("I18").Value = 792

$wsCtrl.Range("I20").Formula = "=SUM(I2:I18)"

$wsCtrl.Range("I27").Select()

# ---------------------------------------------------------------------
# Leave "out_vars" as the active/front sheet, matching the original file.
# ---------------------------------------------------------------------
$wsOut.Activate()
